# Penalty Reward System update
# - "Forecast Comparison" sheet: shift Week_Start_Date forward by one week
#   (each row now shows the date that used to belong to the next row, plus
#   a brand-new date for the final week) and refresh the MyForecast (D)
#   numbers to match.
# - "Summary" sheet: refresh the derived stats (historical range end date,
#   mean/median, totals, max/min forecast + the weeks they occur on) to
#   match the new forecast numbers.
#
# Date-like / number-like text must stay TEXT (matches the source file,
# which stores these as inline strings) -- Excel auto-converts a bare
# "2025-01-12"-style value typed into .Value to a real date serial, so a
# leading apostrophe is used to force literal text entry, same as typing
# '2025-01-12 into a cell in the UI.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Forecast Comparison
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$weekStarts = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

$myForecast = @(75, 75, 74, 76, 78, 80, 77, 71, 69, 73, 78, 81, 79, 75, 73, 77)

for ($i = 0; $i -lt $weekStarts.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = "'" + $weekStarts[$i]
    $ws1.Cells.Item($row, 4).Value = $myForecast[$i]
}

# ---------------------------------------------------------------------
# Sheet 2: Summary
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B2").Value = "2023-02-19 to 2025-01-05"
$ws2.Range("B5").Value = "'49"
$ws2.Range("B6").Value = "'47"
$ws2.Range("B8").Value = "3593 units"
$ws2.Range("B9").Value = "'1211"
$ws2.Range("B10").Value = "'606"
$ws2.Range("B11").Value = "'300"
$ws2.Range("B12").Value = "'81"
$ws2.Range("B13").Value = "'2025-03-30"
$ws2.Range("B14").Value = "'69"
$ws2.Range("B15").Value = "'2025-03-09"
